# correction in sa algorithm and 746 logs
# Update the Fitness column (C) values for rows 2-179 in the log sheet
# to reflect the corrected SA algorithm run results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C63").Value = 7765
$ws.Range("C64:C69").Value = 7345
$ws.Range("C70:C73").Value = 7343
$ws.Range("C74:C179").Value = 7310
